# "Generate Report for Archive" -- refresh the localization-status report:
#   * flip the in-flight status from "Ready for handoff" to "In Translation"
#     on the Overview roll-up sheet and on each per-locale sheet
#   * re-size the status columns to their new (narrower) auto-fit width

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Status values: Overview!E2 (zh-cn column) / Overview!F2 (de-de column),
# and the "Status" column (C2) on each locale's own detail sheet.
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"

# Narrower status/language columns to match the refreshed report layout.
$overview.Columns.Item(5).ColumnWidth = 12.5   # "zh-cn" column
$overview.Columns.Item(6).ColumnWidth = 12.5   # "de-de" column
$zhcn.Columns.Item(3).ColumnWidth = 12.5        # "Status" column
$dede.Columns.Item(3).ColumnWidth = 12.5        # "Status" column
